$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Scenario ID update for the employment by MSA SR14 preliminary forecast:
# year 2016 rows: 75 -> 89
$ws.Range("B2:B8").Value = 89

# year 2020 rows: 82 -> 101
$ws.Range("B9:B15").Value = 101

# year 2025 rows: 81 -> 102
$ws.Range("B16:B22").Value = 102

# year 2035 rows: 76 -> 104
$ws.Range("B23:B29").Value = 104

# Move/restore the active selection on the sheet to D26
$ws.Range("D26").Select()
